$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert a new "2022-Q1" sheet right before the "总计" (summary) sheet.
#    It holds the same per-fund-holding layout as the other quarter
#    sheets (2021-Q2 / 2021-Q3 / 2021-Q4), so start from a copy of the
#    most recent quarter ("2021-Q4") to inherit its header row, fund
#    list and formatting, then patch in the 2022-Q1 numbers.
# ------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($wb.Worksheets.Item("总计"))
$q1.Name = "2022-Q1"

$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Range("A1:H3").Copy($q1.Range("A1:H3"))

# Keep the numeric-looking text columns as literal text (percentages /
# amounts are stored as text in this workbook, same as the other
# quarter sheets) while swapping in this quarter's figures.
$q1.Range("D2:G3").NumberFormat = "@"

$q1.Range("D2").Value = "0.06"
$q1.Range("E2").Value = "94.34"
$q1.Range("F2").Value = "3.23"
$q1.Range("G2").Value = "0.0019"
$q1.Range("H2").Value = 2

$q1.Range("D3").Value = "0.01"
$q1.Range("E3").Value = "94.34"
$q1.Range("F3").Value = "3.23"
$q1.Range("G3").Value = "0.0003"
$q1.Range("H3").Value = 2

# ------------------------------------------------------------------
# 2. Add a matching row at the top of "总计", pushing the existing
#    quarters down by one (2021-Q4/Q3/Q2 each shift down a row).
#    Re-fetch "总计" by name now that the sheet collection has
#    shifted — a reference captured before Add() would otherwise keep
#    tracking the slot the new sheet just took over.
# ------------------------------------------------------------------
$zj = $wb.Worksheets.Item("总计")

$zj.Range("A4:D4").Copy($zj.Range("A5:D5"))
$zj.Range("A3:D3").Copy($zj.Range("A4:D4"))
$zj.Range("A2:D2").Copy($zj.Range("A3:D3"))

$zj.Range("A3").Value = 1
$zj.Range("A4").Value = 2
$zj.Range("A5").Value = 3

$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 0
